$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @(
    @(15, 1, 29.09,  435.548,              1.97736638, 20564362.12, 8.29189352,  15.78603854, 461.9314859,  4.549648183892288, 68.09076589965487,  0.3042690013133379, 7747045.444773559,  1.576608731235504, 2.399488369536128, 116.2805947480786, 0.15),
    @(15, 1, 43.16,  632.6799999999999,    1.3564264,  25707741.752, 11.63860912, 32.91796542, 1428.68794184, 6.669885662400743, 93.46923877508787,  0.1885996524334975, 8483622.617616583,  3.003980909877465, 7.576114327711927, 433.9464282085918, 0.85),
    @(20, 1, 21.694, 432.682,              1.99801192, 34295968.112, 10.0291708,  10.69102954, 234.3309088,  3.624726373474987, 72.16139834303132,  0.3335295188101783, 10749485.79427716,  3.201314198864626, 3.272110775208575, 89.66405358001231, 0.15),
    @(20, 1, 33.416, 645.302,              1.33802164, 48130139.712, 11.79523276, 19.30532634, 650.48072802, 5.830517461403874, 106.7509171512135,  0.2153220606796737, 15414344.87170735,  3.448034262654486, 5.232150290690965, 229.6351214302819, 0.85)
)

$startRow = 14
for ($i = 0; $i -lt $data.Count; $i++) {
    $rowIndex = $startRow + $i
    $rowValues = $data[$i]
    for ($col = 1; $col -le $rowValues.Count; $col++) {
        $ws.Cells.Item($rowIndex, $col).Value = $rowValues[$col - 1]
    }
}

$wb.Save()
